# Updated symbol list (cryptocurrency prices) to match the Dec 30 2022 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text; prefix with a single-quote so Excel keeps
# storing the numeric-looking strings as text (matching the original inlineStr cells)
# instead of converting them to numbers.
$ws.Range("D2").Value = "'243.87"
$ws.Range("D3").Value = "'23.98"
$ws.Range("D4").Value = "'5.128"
$ws.Range("D5").Value = "'0.05749"
$ws.Range("D6").Value = "'6.487"
$ws.Range("D7").Value = "'3.157"
$ws.Range("D8").Value = "'0.8102"
$ws.Range("D9").Value = "'0.8373"
$ws.Range("D10").Value = "'0.1345"
$ws.Range("D11").Value = "'0.06959"
$ws.Range("D12").Value = "'0.03122"
$ws.Range("D13").Value = "'0.02834"
$ws.Range("D14").Value = "'0.09369"
$ws.Range("D16").Value = "'0.001528"
$ws.Range("D17").Value = "'0.04669"
$ws.Range("D18").Value = "'0.0005968"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006081"
$ws.Range("D20").Value = "'0.001240"
$ws.Range("D22").Value = "'0.00008712"
$ws.Range("D23").Value = "'3.501"
$ws.Range("D24").Value = "'2.083"
$ws.Range("D26").Value = "'0.1339"
$ws.Range("D28").Value = "'0.0002332"
$ws.Range("D40").Value = "'0.03626"
$ws.Range("D41").Value = "'0.006308"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1046"
$ws.Range("D43").Value = "'0.002936"
$ws.Range("D44").Value = "'0.007365"
$ws.Range("D45").Value = "'0.00005311"
$ws.Range("D47").Value = "'0.2734"
$ws.Range("D48").Value = "'0.002275"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D50").Value = "'0.0002003"
